# Update column G ("K") values on the active sheet, rows 2-41, per the
# regenerated save_data (K computed in place of the old Strike# based value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 3
    4  = 2
    5  = 0
    6  = 0
    7  = 3
    8  = 2
    9  = 2
    10 = 3
    11 = 1
    12 = 0
    13 = 2
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 2
    20 = 3
    21 = 2
    22 = 3
    23 = 2
    24 = 4
    25 = 3
    26 = 2
    27 = 5
    28 = 2
    29 = 0
    30 = 1
    31 = 4
    32 = 3
    33 = 6
    34 = 4
    35 = 6
    36 = 9
    37 = 6
    38 = 3
    39 = 3
    40 = 2
    41 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
